# Teleop - Fixed elevator WITH Buttons
#
# Merge the run sequence that follows the "rd" superscript run in the
# "On January 3rd, 2015, ..." paragraph (slide 1, TextBox 5) back into a
# single run, matching PowerPoint's normal behaviour after an in-place
# text edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(109)
$tr = $shape.TextFrame.TextRange

$full = $tr.Text

# Locate the end of the "rd" (superscript ordinal) run so we only touch
# the text that comes after it.
$rdIndex = $full.IndexOf("rd")
$afterStart = $rdIndex + 3
$afterLength = $full.Length - $afterStart + 1

$afterRange = $tr.Characters($afterStart, $afterLength)

$newText = ", 2015, our team competed in a New Jersey Qualifier at Liberty Science Center. We used our newly made robot that included a  carbon-fiber chasse, a swerve-drive, and scissor lift that also was not entirely operational. We did emerge successful at the end of the Qualifier, due to strategic driving in tele-op."

$afterRange.Text = $newText
